$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values: TEMP_EXT_HIGHER_LIMIT (C2) and TEMP_EXT_LOWER_LIMIT (D2)
$ws.Range("C2").Value = 28
$ws.Range("D2").Value = 10

# Move the active selection from G3 to C3 (as recorded in the saved view state)
$ws.Range("C3").Select()
